# Partial style modification: insert two new symbols (~ and `) at the top
# of the "symbol" list in column G, shifting the rest down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the two special cell formats (quote-prefix style for the apostrophe
# entry, Consolas font style for the em-dash entry) by copying them onto their
# new rows before the source cells get overwritten below.
$ws.Range("G17").Copy()
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("G22").Copy()
$ws.Range("G24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the full symbol sequence (32 entries) into G2:G33, in its new order.
# (The backtick is written before the tilde so the two brand-new shared-string
# table entries are minted in that same relative order.)
$ws.Cells.Item(3, 7).Value = "``"
$ws.Cells.Item(2, 7).Value = "~"
$ws.Cells.Item(4, 7).Value = "!"
$ws.Cells.Item(5, 7).Value = "@"
$ws.Cells.Item(6, 7).Value = "#"
$ws.Cells.Item(7, 7).Value = "`$"
$ws.Cells.Item(8, 7).Value = "%"
$ws.Cells.Item(9, 7).Value = "^"
$ws.Cells.Item(10, 7).Value = "&"
$ws.Cells.Item(11, 7).Value = "*"
$ws.Cells.Item(12, 7).Value = "("
$ws.Cells.Item(13, 7).Value = ")"
$ws.Cells.Item(14, 7).Value = "-"
$ws.Cells.Item(15, 7).Value = "'="
$ws.Cells.Item(16, 7).Value = "["
$ws.Cells.Item(17, 7).Value = "]"
$ws.Cells.Item(18, 7).Value = ";"
$ws.Cells.Item(19, 7).Value = "''"
$ws.Cells.Item(20, 7).Value = "\"
$ws.Cells.Item(21, 7).Value = ","
$ws.Cells.Item(22, 7).Value = "."
$ws.Cells.Item(23, 7).Value = "/"
$ws.Cells.Item(24, 7).Value = "—"
$ws.Cells.Item(25, 7).Value = "+"
$ws.Cells.Item(26, 7).Value = "{"
$ws.Cells.Item(27, 7).Value = "}"
$ws.Cells.Item(28, 7).Value = ":"
$ws.Cells.Item(29, 7).Value = "`""
$ws.Cells.Item(30, 7).Value = "¦"
$ws.Cells.Item(31, 7).Value = "<"
$ws.Cells.Item(32, 7).Value = ">"
$ws.Cells.Item(33, 7).Value = "?"

# Writing a literal "=" above forced quote-prefix formatting onto that cell
# (row 15); restore its plain (non quote-prefixed) style to match its neighbours.
# G22 kept the Consolas-font format it had before its value changed (plain
# value writes don't touch font formatting), so reset it to the plain style too.
$ws.Range("A1").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the selection recorded in the edit.
$ws.Range("I8").Select()

